# Supplementary_table_template.xlsx edit
#
# Commit: "removed afo, cheminf otw, changed references in pid article,
# updated table in LbE, further minor fixes in mcrs and ctfiles articles"
#
# This removes the AFO (Allotrope Foundation Ontology) related "local
# identifier" / "sample identifier" row entries (columns B and C in rows
# 2-4), which also drops the now-unused AFR_0001118 / AFR_0000919 shared
# strings once Excel recompacts the shared-string table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "local identifier" / "sample identifier" (AFO) columns' data
# from the header block (rows 2-4, columns B and C). Use Clear() (not just
# ClearContents) so the cell's style is dropped too and the <c> element
# disappears entirely, matching the target workbook.
$ws.Range("B2:C4").Clear()

# The author ended up with cell A4 selected when the file was last saved.
$ws.Range("A4").Select()
